$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$about = "Eicher 333 Super Plus (five Star) is an amazing and powerful tractor with a super attractive design. Eicher 333 Super Plus (five Star) is an effective tractor launched by the Eicher Tractor. The 333 Super Plus (five Star) comes with all the advanced technology for effective work on the farm. Here we show all the features, quality, and fair price of the Eicher 333 Super Plus (five Star) Tractor. Check down below.`nEicher 333 Super Plus (five Star) Engine Capacity`nThe tractor comes with 36 HP. Eicher 333 Super Plus (five Star) engine capacity provides efficient mileage on the field. The Eicher 333 Super Plus (five Star) is one of the powerful tractors and offers good mileage. The 333 Super Plus (five Star) Tractor has a capability to provide high performance on the field. Eicher 333 Super Plus (five Star) comes with super power which is fuel efficient.`nEicher 333 Super Plus (five Star) Quality Features`nIt has 8 Forward + 2 Reverse gearboxes.`nAlong with this, Eicher 333 Super Plus (five Star) has a superb kmph forward speed.`nEicher 333 Super Plus (five Star) manufactured with Oil Immersed Brake.`nEicher 333 Super Plus (five Star) Steering type is smooth .`nIt offers a litre large fuel tank capacity for long hours on farms.`nEicher 333 Super Plus (five Star) has 1650 kg strong Lifting capacity.`nThis 333 Super Plus (five Star) tractor consists of multiple tread pattern tyres for the effective work`nEicher 333 Super Plus (five Star) Tractor Price`nEicher 333 Super Plus (five Star) Price in India is Rs. 6.10-6.20 Lakh*. The 333 Super Plus (five Star) price is set according to the budget of Indian farmers. It is the main reason Eicher 333 Super Plus (five Star) became popular among Indian farmers with its launch. For other inquiries related to Eicher 333 Super Plus (five Star), stay tuned with TractorJunction. You can find videos related to the 333 Super Plus (five Star) tractor from which you can get more information about Eicher 333 Super Plus (five Star). Here you can also get an updated Eicher 333 Super Plus (five Star) Tractor on road price 2024.`nWhy Tractor Junction for Eicher 333 Super Plus (five Star)?`nYou can get Eicher 333 Super Plus (five Star) at Tractor Junction with exclusive features. If you have any further queries related to Eicher 333 Super Plus (five Star), you can contact us. Our customer executive will help you out and tell you all about Eicher 333 Super Plus (five Star). So, visit Tractor Junction and get Eicher 333 Super Plus (five Star) with price and features. You can also compare Eicher 333 Super Plus (five Star) with other tractors.`nGet latest Eicher 333 Super Plus (five Star) on road price Jan 20, 2024."

# Update data row (row 2) values to describe the Eicher 333 Super Plus (five Star) tractor instead
# of the previous Sonalika MM+ 39 DI entry. Columns whose value is identical between the old and
# new tractor (C, F, I, M, S, AJ, AP) are intentionally left untouched.
$ws.Range("A2").Value = "Eicher Tractors"
$ws.Range("B2").Value = "333 Super Plus (five Star)"
$ws.Range("D2").Value = "36 HP"
$ws.Range("E2").ClearContents()
$ws.Range("G2").Value = "Oil Immersed Brake"
$ws.Range("H2").Value = "N/A"
$ws.Range("J2").Value = $about
$ws.Range("K2").Value = "2365 CC"
$ws.Range("L2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("Q2").Value = "Centre Shift Option; Side Shift Partial Constant Mesh"
$ws.Range("R2").Value = "Single/Dual"
$ws.Range("V2").Value = "30.84 kmph"
$ws.Range("X2").ClearContents()
$ws.Range("Z2").ClearContents()
$ws.Range("AA2").ClearContents()
$ws.Range("AC2").ClearContents()
$ws.Range("AH2").Value = "1650 kg"
$ws.Range("AI2").ClearContents()
$ws.Range("AK2").ClearContents()
$ws.Range("AL2").ClearContents()
$ws.Range("AM2").ClearContents()
$ws.Range("AO2").ClearContents()

# Setting the long multi-line About text causes Excel to auto-adjust row 2's
# height (customHeight). Re-running AutoFit restores the default (non-custom)
# row height so the saved sheet matches the original formatting.
$ws.Rows.Item(2).AutoFit()
